# fix(gui) step 1 and 2
# Step 1: the four "Almohadilla" autoadhesiva prices go from 167 -> 400
# Step 2: the "ZOCALO" autoadhesivo price goes from 508 -> 1182

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D29").Value = 400
$ws.Range("D30").Value = 400
$ws.Range("D31").Value = 400
$ws.Range("D32").Value = 400

$ws.Range("D33").Value = 1182
